$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "score" column (C) for the three players with their results.
$ws.Range("C2").Style = "Normal"
$ws.Range("C2").Value = "8 / 12"

$ws.Range("C8").Style = "Normal"
$ws.Range("C8").Value = "12 / 12"

$ws.Range("C17").Style = "Normal"
$ws.Range("C17").Value = "2 / 12"

# Mark every logo attempt (column D) as a success (green), except the rows
# handled individually below.
#
# Order matters here: resetting the style to "Normal" before touching the
# font means each genuinely new font color creates a brand new font/style
# entry instead of silently being folded into a pre-existing, coincidentally
# identical one. The first new color used becomes font/style id 20/25, the
# second 21/26, the third 22/27 - so we apply green first, then red, then
# orange, to land on the expected ids.
$successRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,19,20,21,22,23,24)
foreach ($r in $successRows) {
    $cell = $ws.Range("D" + $r)
    $cell.Style = "Normal"
    $cell.Value = "Succès"
    $cell.Font.Color = 32768
}

# Row 27 keeps its original failure message, just recolored red.
$ws.Range("D27").Style = "Normal"
$ws.Range("D27").Value = "Échec : Brand « heineken » is prohibited."
$ws.Range("D27").Font.Color = 255

# The remaining rows get the orange "warning" tone.
$ws.Range("D18").Style = "Normal"
$ws.Range("D18").Value = "Alerte : Brand « fake » does not exist."
$ws.Range("D18").Font.Color = 26367

$ws.Range("D25").Style = "Normal"
$ws.Range("D25").Value = "Succès sur cet élément"
$ws.Range("D25").Font.Color = 26367

$ws.Range("D26").Style = "Normal"
$ws.Range("D26").Value = "Succès sur cet élément"
$ws.Range("D26").Font.Color = 26367

$ws.Range("D28").Style = "Normal"
$ws.Range("D28").Value = "Élément non utilisé par le robot suite à une erreur précédente."
$ws.Range("D28").Font.Color = 26367
